$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1095408508766886
$ws.Range("B3").Value = 0.09923718903447988
$ws.Range("H3").Value = 0.2087780399111685
$ws.Range("B4").Value = 0.08293037643502431
$ws.Range("H4").Value = 0.1924712273117129
$ws.Range("B5").Value = 0.08664581060446666
$ws.Range("H5").Value = 0.1961866614811553
$ws.Range("B6").Value = 0.08809867114853734
$ws.Range("H6").Value = 0.1976395220252259
$ws.Range("B7").Value = 0.08783588393753544
$ws.Range("C7").Value = 0.006160670765040378
$ws.Range("D7").Value = 5.923068308793769
$ws.Range("E7").Value = 0.01595706240556712
$ws.Range("F7").Value = 0.07575650494330734
$ws.Range("G7").Value = 0.09991526293176362
$ws.Range("H7").Value = 0.1973767348142241
$ws.Range("B8").Value = 0.0157241038097211
$ws.Range("C8").Value = 0.001721439986165525
$ws.Range("D8").Value = 1.250819828511136
$ws.Range("E8").Value = 0.01094147031537767
$ws.Range("F8").Value = 0.01234771130724778
$ws.Range("G8").Value = 0.01910049631219455
$ws.Range("H8").Value = 0.1252649546864097
$ws.Range("B9").Value = 0.01208892776099924
$ws.Range("C9").Value = 0.001840075054147123
$ws.Range("D9").Value = 0.87551099422637
$ws.Range("E9").Value = 0.009477498262888963
$ws.Range("F9").Value = 0.008479945585292482
$ws.Range("G9").Value = 0.01569790993670598
$ws.Range("H9").Value = 0.1216297786376879
$ws.Range("B10").Value = 0.01380883928723627
$ws.Range("C10").Value = 0.001659177192505468
$ws.Range("D10").Value = 0.8671545001381232
$ws.Range("E10").Value = 0.01561553626919321
$ws.Range("F10").Value = 0.01055572963197893
$ws.Range("G10").Value = 0.01706194894249302
$ws.Range("H10").Value = 0.1233496901639249
$ws.Range("B11").Value = 0.03167191081166314
$ws.Range("H11").Value = 0.1412127616883518
$ws.Range("B12").Value = 0.05210666110496477
$ws.Range("H12").Value = 0.1616475119816534
$ws.Range("B13").Value = 0.06291557863752552
$ws.Range("H13").Value = 0.1724564295142141
$ws.Range("B14").Value = 0.07365463933776871
$ws.Range("H14").Value = 0.1831954902144573
$ws.Range("B15").Value = 0.07972835105340612
$ws.Range("H15").Value = 0.1892692019300947
$ws.Range("B16").Value = 0.08101929742172134
$ws.Range("H16").Value = 0.1905601482984099
$ws.Range("B17").Value = 0.08334608047783983
$ws.Range("H17").Value = 0.1928869313545284
$ws.Range("B18").Value = -0.1095408508766886
$ws.Range("B19").Value = 0.08588987266372432
$ws.Range("H19").Value = 0.1954307235404129
$ws.Range("B20").Value = 0.08952242401224325
$ws.Range("H20").Value = 0.1990632748889319
$ws.Range("B21").Value = 0.09089316000939306
$ws.Range("H21").Value = 0.2004340108860817
$ws.Range("B22").Value = 0.09475971894242061
$ws.Range("C22").Value = 0.007620871415397229
$ws.Range("D22").Value = 1009979561198.836
$ws.Range("E22").Value = 0.04347295510024996
$ws.Range("F22").Value = 0.07975882774137004
$ws.Range("G22").Value = 0.1097606101434712
$ws.Range("H22").Value = 0.2043005698191092
$ws.Range("B23").Value = 0.09704936308425112
$ws.Range("H23").Value = 0.2065902139609397
$ws.Range("B24").Value = 0.1001725386056939
$ws.Range("C24").Value = 0.007814031344815753
$ws.Range("D24").Value = 996697712547.6837
$ws.Range("E24").Value = 0.05782375755067647
$ws.Range("F24").Value = 0.08481222431341105
$ws.Range("G24").Value = 0.1155328528979767
$ws.Range("H24").Value = 0.2097133894823825
$ws.Range("B25").Value = 0.1021566949567594
$ws.Range("C25").Value = 0.007892856342740262
$ws.Range("D25").Value = 603442959324.0503
$ws.Range("E25").Value = 0.05376238981176969
$ws.Range("F25").Value = 0.08662495127093868
$ws.Range("G25").Value = 0.1176884386425802
$ws.Range("H25").Value = 0.2116975458334481
$ws.Range("B26").Value = 0.1012313291632006
$ws.Range("C26").Value = 0.008015375722886199
$ws.Range("D26").Value = 22.82675888221823
$ws.Range("E26").Value = 0.05990244578407614
$ws.Range("F26").Value = 0.08547361794743422
$ws.Range("G26").Value = 0.116989040378967
$ws.Range("H26").Value = 0.2107721800398892
$ws.Range("B27").Value = 0.1076876494503328
$ws.Range("C27").Value = 0.008049226972368427
$ws.Range("D27").Value = 22.76106326311974
$ws.Range("E27").Value = 0.06557252171094248
$ws.Range("F27").Value = 0.09185393207197737
$ws.Range("G27").Value = 0.1235213668286879
$ws.Range("H27").Value = 0.2172285003270214
$ws.Range("B28").Value = 0.1076451618059842
$ws.Range("C28").Value = 0.008220632994280091
$ws.Range("D28").Value = 951876673633.943
$ws.Range("E28").Value = 0.09908224531215967
$ws.Range("F28").Value = 0.09149338219573487
$ws.Range("G28").Value = 0.1237969414162337
$ws.Range("H28").Value = 0.2171860126826728
$ws.Range("B29").Value = 0.0177844887245429
$ws.Range("C29").Value = 0.001713823760253659
$ws.Range("D29").Value = 1.717639576708349
$ws.Range("E29").Value = 0.0102958081672072
$ws.Range("F29").Value = 0.01438747950771367
$ws.Range("G29").Value = 0.02118149794137169
$ws.Range("H29").Value = 0.1273253396012315
